$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 84

# Column A holds a date-like label that must stay plain text (matches the
# existing rows, which are text, not real Excel dates), so force a text
# number format before assigning the value to stop Excel's autoconvert.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "02/16/2026"

$ws.Cells.Item($row, 2).Value  = 9608.030000000001
$ws.Cells.Item($row, 3).Value  = 0.2431908283538754
$ws.Cells.Item($row, 4).Value  = 0.7568091716461246
$ws.Cells.Item($row, 5).Value  = -307.54
$ws.Cells.Item($row, 6).Value  = -34.75
$ws.Cells.Item($row, 7).Value  = -23670.31
$ws.Cells.Item($row, 8).Value  = -76.5
$ws.Cells.Item($row, 9).Value  = -1066.26
$ws.Cells.Item($row, 10).Value = -31.33
$ws.Cells.Item($row, 11).Value = -24736.57
$ws.Cells.Item($row, 12).Value = -72.02
